$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P-values and Odds Ratios")

# ---------------------------------------------------------------------------
# Row 79: stray note cell (set first so new shared-string order matches
# the source workbook: ">" , "19 or less", "60 or more")
# ---------------------------------------------------------------------------
$ws.Range("B79").Value = ">"

# ---------------------------------------------------------------------------
# Row 42: section header labels ("Without Naxalone" / "With Naxalone")
# (copy formats from the equivalent B16:C16 header cells above so B42/C42
# keep the same placeholder style as the rest of the sheet)
# ---------------------------------------------------------------------------
$ws.Range("B16:C16").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("D42").Value = "Without Naxalone"
$ws.Range("F42").Value = "With Naxalone"

# ---------------------------------------------------------------------------
# Row 43: column headers (bold, reuse existing style seen in rows 17/23/30/36)
# ---------------------------------------------------------------------------
$ws.Range("C43").Value = "# Victims"
$ws.Range("D43").Value = "# Died"
$ws.Range("E43").Value = "% Died"
$ws.Range("F43").Value = "# Victims"
$ws.Range("G43").Value = "# Died"
$ws.Range("H43").Value = "% Died"
$ws.Range("I43").Value = "Unadj Odds Ratio"
$ws.Range("C43:I43").Font.Bold = $true
$ws.Range("J43").Value = "Odds Ratio"
$ws.Range("K43").Value = "Lower 95% CI"
$ws.Range("L43").Value = "Upper 95% CI"

# ---------------------------------------------------------------------------
# Rows 44-50: Age-group breakdown data (victims/deaths with/without Naloxone)
# ---------------------------------------------------------------------------
$ageLabels = @("19 or less", "20 - 24", "25 - 29", "30 - 39", "40 - 49", "50 - 59", "60 or more")
$noNarcanVictims = @(96, 348, 606, 989, 467, 269, 89)
$noNarcanDied    = @(10, 105, 190, 354, 203, 131, 45)
$narcanVictims   = @(104, 724, 1338, 2109, 793, 452, 141)
$narcanDied      = @(9, 30, 76, 148, 69, 28, 7)

for ($i = 0; $i -lt 7; $i++) {
    $r = 44 + $i

    $ws.Range("B$r").Value = $ageLabels[$i]

    $ws.Range("C$r").Value = $noNarcanVictims[$i]
    $ws.Range("D$r").Value = $noNarcanDied[$i]
    $ws.Range("E$r").Formula = "=D$r/C$r"
    $ws.Range("E$r").NumberFormat = "0.00%"

    $ws.Range("F$r").Value = $narcanVictims[$i]
    $ws.Range("G$r").Value = $narcanDied[$i]
    $ws.Range("H$r").Formula = "=G$r/F$r"
    $ws.Range("H$r").NumberFormat = "0.00%"

    $ws.Range("I$r").Formula = '=CONCATENATE(TEXT(J' + $r + ',"0.00")," (",TEXT(K' + $r + ',"0.00"),"-",TEXT(L' + $r + ',"0.00"),")")'

    $ws.Range("J$r").Formula = "=H$r/(1-H$r)*(1-E$r)/E$r"
    $ws.Range("J$r").NumberFormat = "0.00"

    $ws.Range("K$r").Formula = "=EXP(LN(J$r)-1.96*SQRT((1/D$r+1/(C$r-D$r)+1/G$r+1/(F$r-G$r))))"
    $ws.Range("K$r").NumberFormat = "0.00"

    $ws.Range("L$r").Formula = "=EXP(LN(J$r)+1.96*SQRT(1/D$r+1/(C$r-D$r)+1/G$r+1/(F$r-G$r)))"
    $ws.Range("L$r").NumberFormat = "0.00"
}

# ---------------------------------------------------------------------------
# Restore the view state (selection) similar to the saved workbook
# ---------------------------------------------------------------------------
$ws.Range("B53").Select()
